$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 3
$ws.Range("B8").Value = 3
$ws.Range("B9").Value = 3
$ws.Range("B10").Value = 3
$ws.Range("B11").Value = 4
$ws.Range("B12").Value = 4
